# Sheet to Object[][] in main methods
# Append additional invalid-login test rows (peter124 .. peter145) to the
# "invalidLoginTest" sheet (Sheet2), then select column A as the last UI action.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$startUser = 124
$endUser = 145
$startRow = 4

for ($n = $startUser; $n -le $endUser; $n++) {
    $row = $startRow + ($n - $startUser)
    $ws.Cells.Item($row, 1).Value = "Peter"
    $ws.Cells.Item($row, 2).Value = "peter$n"
    $ws.Cells.Item($row, 3).Value = "Invalid credentials"
}

$ws.Activate()
$ws.Columns("A").Select()
